$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "Lettuce-Refrigerator,tomato-Cabinet,parmesan cheese-Refrigerator,caesar dressing-Refrigerator, chicken breast-Refrigerator, croutons-Cabinet"
$ws.Range("C2").Value = "Salmon-Refrigerator,lettuce-Refrigerator,cucumber-Refrigerator,tomato-Cabinet,feta cheese-Refrigerator, red onion-Refrigerator, sweet green pepper-Refrigerator, chickpeas-Cabinet"
$ws.Range("C3").Value = "Chicken breast-Refrigerator,basmati rice-Cabinet, green beans-Freezer, chipotle seasoning-Cabinet"
$ws.Range("C4").Value = "Bison filet-Refrigerator,salt-Cabinet,pepper-Cabinet,lettuce-Refrigerator,walnuts-Cabinet,strawberries-Refrgierator,balsamic dressing-Refrigerator "
$ws.Range("C5").Value = "Salmon-Refrigerator, salt-Cabinet, pepper-Cabinet, quinoa-Cabinet, asparagus-Refrigerator"
$ws.Range("C6").Value = "Crab cakes-Refrigerator,russett potato-Cabinet, butter-Refrigerator, salt-Cabinet, pepper-Cabinet,brussel sprouts-Refrigerator"
$ws.Range("C7").Value = "Maple syrup-Cabinet, soy sauce-Refrigerator, dijon mustard-Refrigerator, garlic-Cabinet, Salmon-Refrigerator, quinoa-Cabinet, peas-Freezer"
